# Update "想去人数" (F column) figures across sheets, matching the
# gh-pages data refresh captured in commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 104
$ws.Range("F4").Value = 415
$ws.Range("F5").Value = 191
$ws.Range("F6").Value = 130
$ws.Range("F7").Value = 1126
$ws.Range("F8").Value = 377
$ws.Range("F9").Value = 190
$ws.Range("F12").Value = 372
$ws.Range("F13").Value = 383
$ws.Range("F14").Value = 782
$ws.Range("F15").Value = 167
$ws.Range("F17").Value = 277
$ws.Range("F18").Value = 74
$ws.Range("F19").Value = 999
$ws.Range("F20").Value = 453
$ws.Range("F21").Value = 259
$ws.Range("F23").Value = 377
$ws.Range("F24").Value = 27
$ws.Range("F26").Value = 466

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 363
$ws.Range("F7").Value = 281
$ws.Range("F11").Value = 149

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 104
$ws.Range("F6").Value = 415
$ws.Range("F7").Value = 191
$ws.Range("F8").Value = 130
$ws.Range("F9").Value = 1126
$ws.Range("F10").Value = 377
$ws.Range("F11").Value = 190
$ws.Range("F14").Value = 363
$ws.Range("F17").Value = 372
$ws.Range("F19").Value = 281
$ws.Range("F20").Value = 383
$ws.Range("F21").Value = 782
$ws.Range("F22").Value = 167
$ws.Range("F24").Value = 277
$ws.Range("F25").Value = 74
$ws.Range("F26").Value = 999
$ws.Range("F27").Value = 453
$ws.Range("F30").Value = 259
$ws.Range("F32").Value = 377
$ws.Range("F34").Value = 149
$ws.Range("F35").Value = 27
$ws.Range("F38").Value = 466
